{"js": "// Portuguese translation pass for the partner \"flight & accommodation details\"\n// email template. Every change below is a like-for-like text swap inside an\n// existing run (formatting/highlights are left untouched because we replace\n// the run's own range, not the whole paragraph where other runs must survive).\n\nconst body = context.document.body;\n\n// Small helper: search for an exact (case-sensitive, literal) run of text and\n// replace the Nth match (0-based, in document order) with newText.\nasync function replaceNth(searchText, index, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWildcards: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length <= index) {\n    throw new Error(`Expected match #${index} for ${JSON.stringify(searchText)}, found ${results.items.length}`);\n  }\n  results.items[index].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// Same, but replaces every match (used when all occurrences need the exact\n// same replacement, e.g. \" or \" -> \" ou \" which appears twice).\nasync function replaceAll(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWildcards: false });\n  results.load(\"text\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// 1. The standalone \"English\" heading right above the table (NOT the\n//    \"English\" inside the language-picker hyperlink at the very top).\nawait replaceNth(\"English\", 1, \"Ingl\u00eas\");\n\n// 2. \"We can't wait to meet you! \" heading.\nawait replaceNth(\"We can\\u2019t wait to meet you! \", 0, \"Mal podemos esperar para estar consigo! \");\n\n// 3. \"Hi \" / \"[PARTNER NAME]\" greeting line.\nawait replaceNth(\"Hi \", 0, \"Ol\u00e1 \");\nawait replaceNth(\"[PARTNER NAME]\", 0, \"[NOME DO PARCEIRO]\");\n\n// 4. \"We hope you're as excited...\" + the second \"[EVENT NAME]\" placeholder\n//    (the first one, in the Subject line, must stay in English).\nawait replaceNth(\"We hope you\\u2019re as excited as we are for \", 0, \"Esperamos que esteja t\u00e3o entusiasmado quanto n\u00f3s com a \");\nawait replaceNth(\"[EVENT NAME]\", 1, \"[NOME DO EVENTO]\");\n\n// 5. \"In this email, we've linked/attached the following documents:\"\nawait replaceNth(\n  \"In this email, we\\u2019ve linked/attached the following documents:\",\n  0,\n  \"Neste e-mail, anex\u00e1mos os seguintes documentos:\"\n);\n\n// 6. Bulleted list items.\nawait replaceNth(\"Your return flight tickets\", 0, \"Os seus bilhetes de avi\u00e3o de regresso\");\nawait replaceNth(\"Your accommodation booking details\", 0, \"Dados da sua reserva de alojamento\");\n\n// 7. \"If you have any questions, please contact us via live chat or WhatsApp.\"\nawait replaceNth(\n  \"If you have any questions, please contact us via \",\n  0,\n  \"Para mais informa\u00e7\u00f5es, contacte-nos atrav\u00e9s de \"\n);\n\n// 8. \"If you have any questions, please contact your country manager, ...\"\nawait replaceNth(\n  \"If you have any questions, please contact your country manager, \",\n  0,\n  \"Para mais quest\u00f5es, pode tamb\u00e9m contactar o seus gestor de parcerias \"\n);\nawait replaceNth(\", at \", 0, \", em \");\n\n// Both remaining \" or \" connectors (live chat/WhatsApp line, and the\n// email/WhatsApp number line) become \" ou \".\nawait replaceAll(\" or \", \" ou \");\n\n// 9. Sign-off line.\nawait replaceNth(\"See you on the \", 0, \"Vemo-nos dia \");\n", "ps1": "# Portuguese translation pass for the partner \"flight & accommodation details\"\n# email template. Every change is a like-for-like text swap done with\n# Find/Replace scoped to a single paragraph's Range, so the two \"English\"\n# strings (one inside the top hyperlink, one as a standalone heading) and the\n# two \"[EVENT NAME]\" placeholders (Subject line vs. body) are not conflated,\n# and so only the exact substring is touched (other runs / formatting, e.g.\n# the yellow highlight on bracket placeholders, are left alone).\n\n$d = $word.ActiveDocument\n\nfunction Replace-InParagraph {\n    param(\n        [int]$ParaIndex,\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n    $rng = $d.Paragraphs.Item($ParaIndex).Range\n    # wdFindStop (0) keeps the search confined to $rng instead of wrapping\n    # around the whole document; wdReplaceAll (2) with MatchCase = $true and\n    # MatchWholeWord = $false (several targets are multi-word phrases with a\n    # leading/trailing space, e.g. \" or \", which can never be a \"whole word\").\n    $ok = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 0, $false, $ReplaceText, 2)\n    if (-not $ok) {\n        throw \"Replace-InParagraph: could not find '$FindText' in paragraph $ParaIndex\"\n    }\n}\n\n# 1. The standalone \"English\" heading right above the table (paragraph 3).\n#    The \"English\" inside the language-picker hyperlink (paragraph 1) is left\n#    untouched because the search is scoped to paragraph 3 only.\nReplace-InParagraph 3 \"English\" \"Ingl\u00eas\"\n\n# 2. \"We can't wait to meet you! \" heading (paragraph 14).\nReplace-InParagraph 14 \"We can\u2019t wait to meet you! \" \"Mal podemos esperar para estar consigo! \"\n\n# 3. \"Hi [PARTNER NAME], \" greeting line (paragraph 15).\nReplace-InParagraph 15 \"Hi \" \"Ol\u00e1 \"\nReplace-InParagraph 15 \"[PARTNER NAME]\" \"[NOME DO PARCEIRO]\"\n\n# 4. \"We hope you're as excited...\" paragraph (paragraph 17), including its\n#    own \"[EVENT NAME]\" placeholder. The Subject line's \"[EVENT NAME]\"\n#    (paragraph 13) is left in English because it's outside this scope.\nReplace-InParagraph 17 \"We hope you\u2019re as excited as we are for \" \"Esperamos que esteja t\u00e3o entusiasmado quanto n\u00f3s com a \"\nReplace-InParagraph 17 \"[EVENT NAME]\" \"[NOME DO EVENTO]\"\n\n# 5. \"In this email, we've linked/attached the following documents:\" (paragraph 19).\nReplace-InParagraph 19 \"In this email, we\u2019ve linked/attached the following documents:\" \"Neste e-mail, anex\u00e1mos os seguintes documentos:\"\n\n# 6. Bulleted list items (paragraphs 20-21).\nReplace-InParagraph 20 \"Your return flight tickets\" \"Os seus bilhetes de avi\u00e3o de regresso\"\nReplace-InParagraph 21 \"Your accommodation booking details\" \"Dados da sua reserva de alojamento\"\n\n# 7. \"If you have any questions, please contact us via live chat or WhatsApp.\" (paragraph 24).\nReplace-InParagraph 24 \"If you have any questions, please contact us via \" \"Para mais informa\u00e7\u00f5es, contacte-nos atrav\u00e9s de \"\nReplace-InParagraph 24 \" or \" \" ou \"\n\n# 8. \"If you have any questions, please contact your country manager, ...\" (paragraph 25).\nReplace-InParagraph 25 \"If you have any questions, please contact your country manager, \" \"Para mais quest\u00f5es, pode tamb\u00e9m contactar o seus gestor de parcerias \"\nReplace-InParagraph 25 \", at \" \", em \"\nReplace-InParagraph 25 \" or \" \" ou \"\n\n# 9. Sign-off line (paragraph 26).\nReplace-InParagraph 26 \"See you on the \" \"Vemo-nos dia \"\n"}
